# MAJ liste staff / CO
# Update staff contact info (phone numbers / names) on the "CO" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO")

# Row 7  - Project manager - Tour de la Releve / Antoine St-Jean
$ws.Range("D7").Value = '\(819) 744-7765'

# Row 8  - Reception and credentials / Claudette Aylwin
$ws.Range("D8").Value = '\(819) 856-0772'

# Row 10 - Race Headquarters / Centrale administrative : name removed, back to TBD
$ws.Range("C10").Value = 'TBD'

# Row 11 - Communications / Press : Sophie-Kristine Richard
$ws.Range("D11").Value = '\(418) 690-6564'

# Row 13 - Lodging Team / Hebergement : new team + phone
$ws.Range("C13").Value = 'Équipe Subway<br/>Mathieu Roy'
$ws.Range("D13").Value = '\(819) 444-9746'

# Row 14 - Medical / Medical : Marie-Eve Dionne & Tina Aube
$ws.Range("D14").Value = '\(819) 860-8614<br/>(819) 860-9605'

# Row 16 - Protocol / Protocole : Julie Pelletier
$ws.Range("D16").Value = '\(819) 856-6200'

# Row 17 - VIP Room / Salon VIP : Claude Deraps
$ws.Range("D17").Value = '\(819) 856-8665'

# Row 18 - Security / Securite : Steve Roussil & Caroline Comeau
$ws.Range("D18").Value = '\(819 825-2210<br/>(819) 860-1580'

# Row 19 - Transportation / Transport : renamed from Bobby Noury to Robert Noury
$ws.Range("C19").Value = 'Robert Noury'

# Row 20 - Facilities and logistics / Equipement et logistique : Michel Laverdure
$ws.Range("D20").Value = '\(819) 856-1565'

# Row 21 - City of Val-d'Or Representative : Luc Lavoie
$ws.Range("D21").Value = '\(819) 856-5375'

# Row 22 - Environmentally friendly / Comite Ecoresponsable : Antoine St-Jean
$ws.Range("D22").Value = '\(819) 744-7765'

# Update the active selection on the CO sheet to match the saved view state.
$ws.Activate()
$ws.Range("B11").Select()
